$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Image folder cleanup: the stimulus that used to be at row 2 (img_5jy9c.png)
# now belongs at row 24, and the one that used to be at row 24 (img_8dmpq.png)
# now belongs at row 2 - swap the stimulus name + all associated score columns
# (L:V) between the two rows.

$row2Vals  = $ws.Range("L2:V2").Value2
$row24Vals = $ws.Range("L24:V24").Value2

$ws.Range("L2:V2").Value2   = $row24Vals
$ws.Range("L24:V24").Value2 = $row2Vals

# Also fix the catch-trial filename on row 16 (text label removed from image).
$ws.Range("L16").Value2 = "stimuli/catch_01.jpg"
